# The presentation is already open as $ppt.ActivePresentation
$p = $ppt.ActivePresentation

# Slide 2 ("How many times has Priya done transactions ...") has a
# duplicated text box ("Rectangle 4") carrying the same caption as
# "Rectangle 5", and the "Chart 6" graphic frame needs to shift left.
$s = $p.Slides.Item(2)

# Remove the duplicate "Rectangle 4" shape (shape index 3: Rectangle 2,
# Chart 3, Rectangle 4, Rectangle 5, Chart 6).
$dupe = $s.Shapes.Item("Rectangle 4")
$dupe.Delete()

# Shift the "Chart 6" chart frame to the left (its top/width/height are
# unchanged). 539552 EMU == 42.4844... pt; the literal below is chosen so
# that the points -> EMU round trip lands exactly on 539552.
$chart = $s.Shapes.Item("Chart 6")
$chart.Left = 42.4844493488189
